$wb = $excel.ActiveWorkbook

# Work on the "field-domain" sheet: insert a new "isPrivate" column between
# the existing "isDomain" (C) and "description" (D) columns.
$fieldDomainWs = $wb.Worksheets.Item("field-domain")

# Inserting a whole column at D shifts the old D1 ("description") to E1
# and keeps its formatting (shared-string value + cell style) intact.
$fieldDomainWs.Range("D1").EntireColumn.Insert()

# Give the freshly inserted D1 the same formatting as its neighbour (C1,
# the bold/red header style) before writing the new header text into it.
$fieldDomainWs.Range("C1").Copy()
$fieldDomainWs.Range("D1").PasteSpecial(-4122)
$fieldDomainWs.Range("D1").Value = "isPrivate"

# Make "field-domain" the active/selected sheet with D1 as the selection,
# which also clears tabSelected/active status from whichever sheet ("info")
# used to be active.
$fieldDomainWs.Activate()
$fieldDomainWs.Range("D1").Select()
